$d = $word.ActiveDocument

function Split-AtRange($rng) {
    # Forces Word to split the run at the boundaries of $rng without
    # changing any visible formatting (toggle Bold on then back off).
    $rng.Bold = 1
    $rng.Bold = 0
}

function Replace-WithRuns($findText, $pieces) {
    # $pieces is an array of strings; the found range's text is replaced
    # by concatenating all pieces, then each piece becomes its own run.
    $full = [string]::Join("", $pieces)
    $r = $d.Content
    $r.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $full, 2) | Out-Null
    $start = $r.Start
    $pos = $start
    foreach ($p in $pieces) {
        $len = $p.Length
        $sub = $d.Range($pos, $pos + $len)
        Split-AtRange $sub
        $pos = $pos + $len
    }
}

# ---------------------------------------------------------------------
# Change 1: "upsert (merge) method" -> "upsert (MERGE) method" with MERGE
# split into its own run.
# ---------------------------------------------------------------------
$r1 = $d.Content
$r1.Find.Execute("upsert (merge) method", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$mergeStart = $r1.Start + "upsert (".Length
$mergeEnd = $mergeStart + "merge".Length
$mergeRange = $d.Range($mergeStart, $mergeEnd)
$mergeRange.Text = "MERGE"
Split-AtRange $mergeRange

Write-Output "change1 done"

# ---------------------------------------------------------------------
# Change 2: "With the data successfully loaded into Snowflake, we create
# two separate views..." paragraph rewritten, with SYMBOL / DATE / CLOSE
# split into their own runs.
# ---------------------------------------------------------------------
$oldPara2 = "With the data successfully loaded into Snowflake, we create two separate views based on the stock symbol" + [char]0x2014 + "Apple (AAPL) and Five Below (FIVE). This separation is necessary because combining both stocks in a single forecast would result in predicting the average price of the two, which is more suitable for analyzing overall market trends rather than individual stock performance. While such an approach could be useful for broader industry-level analysis (e.g., comparing multiple tech stocks to assess market trends), our focus is on forecasting individual stock prices."

$pieces2 = @(
    "With the data successfully loaded into Snowflake, the next step is to create a view to structure the dataset for forecasting. In Snowflake, we can leverage the SERIES_COLNAME parameter to automatically separate data based on the stock symbol, allowing us to handle multiple stocks within the same dataset efficiently. Given this capability, our view only needs to isolate key columns" + [char]0x2014,
    "SYMBOL",
    ", ",
    "DATE",
    ", and ",
    "CLOSE",
    [char]0x2014 + "which are essential for performing time-series forecasting."
)
Replace-WithRuns $oldPara2 $pieces2

Write-Output "change2 done"

# ---------------------------------------------------------------------
# Change 3: "Since the stock symbol was stored..." paragraph rewritten,
# with CLOSE / DATE split into their own runs.
# ---------------------------------------------------------------------
$oldPara3 = "Since the stock symbol was stored as a separate column in the DataFrame, we can efficiently filter records using a WHERE clause in SQL to ensure that each view contains data for only one stock. Additionally, Snowflake" + [char]0x2019 + "s forecasting function requires users to specify both a timestamp column and a target variable for analysis. Extracting and structuring these columns properly is essential to ensure the data is optimized for forecasting and predictive analytics."

$pieces3 = @(
    "Once the view is created, we can apply Snowflake" + [char]0x2019 + "s built-in forecasting function to generate stock price predictions. By specifying the target column (",
    "CLOSE",
    ") and timestamp (",
    "DATE",
    "), Snowflake automatically detects the data intervals within the dataset, ensuring accurate time-series analysis. Furthermore, Snowflake applies multiple machine learning models to the data and automatically selects the most effective one based on performance metrics."
)
Replace-WithRuns $oldPara3 $pieces3

Write-Output "change3 done"

# ---------------------------------------------------------------------
# Change 4: delete the whole "After creating the necessary views..."
# paragraph.
# ---------------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("After creating the necessary views", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$r4.Expand(4) | Out-Null
$r4.Delete()

Write-Output "change4 done"

# ---------------------------------------------------------------------
# Change 5: "...we will have one" / " have two separate files..." ->
# "...we will have" / " two separate files..." (keep the two-run split).
# ---------------------------------------------------------------------
$r5 = $d.Content
$r5.Find.Execute("we will have one", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$delA5 = $d.Range($r5.End - 4, $r5.End)
$delA5.Delete()
$runBStart5 = $r5.End - 4
$delB5 = $d.Range($runBStart5 + 1, $runBStart5 + 1 + 5)
$delB5.Delete()

$para5 = $d.Range($runBStart5, $runBStart5)
$para5.Expand(4) | Out-Null
$splitRange5 = $d.Range($runBStart5, $para5.End - 1)
Split-AtRange $splitRange5

Write-Output "change5 done"

# ---------------------------------------------------------------------
# Change 6: remove the <w:lastRenderedPageBreak/> before "executed
# sequentially..." -- merges that run with the preceding one, but the
# following ". " runs must stay split exactly as they were.
# ---------------------------------------------------------------------
$r6 = $d.Content
$old6 = "Since these operations must be executed sequentially, Apache Airflow" + [char]0x2019 + "s bitshift operator (>>) will be used to ensure that model training runs first, followed by the prediction task, which depends on the model" + [char]0x2019 + "s results"
$r6.Find.Execute($old6, $true, $false, $false, $false, $false, $true, 1, $false, $old6, 2) | Out-Null
$dotStart6 = $r6.End

$p6 = $d.Range($dotStart6, $dotStart6)
$p6.Expand(4) | Out-Null

$dotRange6 = $d.Range($dotStart6, $dotStart6 + 1)
Split-AtRange $dotRange6

$spaceRange6 = $d.Range($dotStart6 + 1, $p6.End - 1)
Split-AtRange $spaceRange6

Write-Output "change6 done"
